# feat. Revise Router Position
#
# Insert a new "stockNum" column immediately before the existing "color"
# column (G). Excel's column-insert shifts color/avgRate/imgUrl/"sql
# start"/sql one column to the right (G->H, H->I, I->J, J->K, K->L) and
# auto-updates every dependent formula/reference in place. The new column
# is then populated with its header and a repeating 1-5 stock count.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at G; existing G:L shift right to H:M and every
# formula/reference that pointed at the old columns is rewritten by Excel
# (e.g. CONCATENATE($J$2,...) becomes CONCATENATE($K$2,...), G9 -> H9, etc.)
$ws.Columns("G:G").Insert() | Out-Null

# Header for the newly inserted column
$ws.Range("G1").Value = "stockNum"

# Populate the new stockNum column with a repeating 1..5 pattern
for ($r = 2; $r -le 39; $r++) {
    $ws.Cells.Item($r, 7).Value = ((($r - 2) % 5) + 1)
}

# Move the active selection to match the edited workbook's view state
$ws.Range("G40").Select() | Out-Null
